## Update Tnfsf14-Tnfrsf14.xlsx with newly-computed TPM-based NATMI values.
##
## The original sheet held 6 sending/target-cluster combinations:
##   rows 2-4 : Sending cluster = ECs   (old TPM numbers)
##   rows 5-7 : Sending cluster = FAPs  (old TPM numbers)
##
## The refreshed script only emits the FAPs-sourced combinations (with
## updated TPM-derived statistics), so the ECs rows are dropped and the
## FAPs rows shift up to become rows 2-4, with several statistic columns
## recalculated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old "ECs" sending-cluster rows (2-4); the old "FAPs" rows
# (5-7) shift up to become rows 2-4, already holding the correct
# Sending/Ligand/Receptor/Target cluster labels and several already-
# matching numeric columns.
$ws.Rows("2:4").Delete()

# Recalculated statistic columns for each of the three remaining rows.
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("M2").Value = 3.961421333333333
$ws.Range("N2").Value = 11.884264
$ws.Range("O2").Value = 0.3114993985605504
$ws.Range("P2").Value = 0.3114993985605504
$ws.Range("Q2").Value = 0.2373736481884444
$ws.Range("R2").Value = 2.136362833696
$ws.Range("S2").Value = 0.3114993985605504
$ws.Range("T2").Value = 0.3114993985605504

$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("O3").Value = 0.5009735319462221
$ws.Range("P3").Value = 0.500973531946222
$ws.Range("S3").Value = 0.5009735319462221
$ws.Range("T3").Value = 0.500973531946222

$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("O4").Value = 0.1875270694932276
$ws.Range("P4").Value = 0.1875270694932276
$ws.Range("S4").Value = 0.1875270694932276
$ws.Range("T4").Value = 0.1875270694932276
